$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column "Matières enseignés" for professors (E1 header)
$ws.Range("E1").Value = "Matières enseignés"

# Widen the CIN/DEPARTEMENT/new column area to fit the content
$ws.Columns.Item(3).ColumnWidth = 26.6
$ws.Columns.Item(4).ColumnWidth = 14.85
$ws.Columns.Item(5).ColumnWidth = 30.85

# Leave the selection where the author left it when saving
$ws.Range("E6").Select()
